# Add three new data columns (HC, HD, HE) to the stats sheet,
# extending the match-by-match table that currently ends at column HB.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# HC = column 211, HD = column 212, HE = column 213
$ws.Cells.Item(1, 211).Value = 10246
$ws.Cells.Item(1, 212).Value = 10253
$ws.Cells.Item(1, 213).Value = 10262
$ws.Cells.Item(2, 211).Value = 2020
$ws.Cells.Item(2, 212).Value = 2020
$ws.Cells.Item(2, 213).Value = 2020
$ws.Cells.Item(3, 211).Value = 10
$ws.Cells.Item(3, 212).Value = 11
$ws.Cells.Item(3, 213).Value = 12
$ws.Cells.Item(4, 211).Value = 1
$ws.Cells.Item(4, 212).Value = 0
$ws.Cells.Item(4, 213).Value = 1
$ws.Cells.Item(5, 211).Value = 0
$ws.Cells.Item(5, 212).Value = 0
$ws.Cells.Item(5, 213).Value = 1
$ws.Cells.Item(6, 211).Value = 41
$ws.Cells.Item(6, 212).Value = 96
$ws.Cells.Item(6, 213).Value = 53
$ws.Cells.Item(7, 211).Value = 82
$ws.Cells.Item(7, 212).Value = 72
$ws.Cells.Item(7, 213).Value = 52
$ws.Cells.Item(8, 211).Value = -41
$ws.Cells.Item(8, 212).Value = 24
$ws.Cells.Item(8, 213).Value = 1
$ws.Cells.Item(9, 211).Value = 0
$ws.Cells.Item(9, 212).Value = 1
$ws.Cells.Item(9, 213).Value = 1
$ws.Cells.Item(10, 211).Value = 14
$ws.Cells.Item(10, 212).Value = 18
$ws.Cells.Item(10, 213).Value = 12
$ws.Cells.Item(11, 211).Value = 211
$ws.Cells.Item(11, 212).Value = 182
$ws.Cells.Item(11, 213).Value = 167
$ws.Cells.Item(12, 211).Value = 87
$ws.Cells.Item(12, 212).Value = 108
$ws.Cells.Item(12, 213).Value = 98
$ws.Cells.Item(13, 211).Value = 298
$ws.Cells.Item(13, 212).Value = 290
$ws.Cells.Item(13, 213).Value = 265
$ws.Cells.Item(14, 211).Value = 2.43
$ws.Cells.Item(14, 212).Value = 1.69
$ws.Cells.Item(14, 213).Value = 1.7
$ws.Cells.Item(15, 211).Value = 119
$ws.Cells.Item(15, 212).Value = 72
$ws.Cells.Item(15, 213).Value = 82
$ws.Cells.Item(16, 211).Value = 34
$ws.Cells.Item(16, 212).Value = 50
$ws.Cells.Item(16, 213).Value = 49
$ws.Cells.Item(17, 211).Value = 23
$ws.Cells.Item(17, 212).Value = 42
$ws.Cells.Item(17, 213).Value = 36
$ws.Cells.Item(18, 211).Value = 18
$ws.Cells.Item(18, 212).Value = 9
$ws.Cells.Item(18, 213).Value = 13
$ws.Cells.Item(19, 211).Value = 11
$ws.Cells.Item(19, 212).Value = 12
$ws.Cells.Item(19, 213).Value = 19
$ws.Cells.Item(20, 211).Value = 4
$ws.Cells.Item(20, 212).Value = 14
$ws.Cells.Item(20, 213).Value = 7
$ws.Cells.Item(21, 211).Value = 1
$ws.Cells.Item(21, 212).Value = 9
$ws.Cells.Item(21, 213).Value = 5
$ws.Cells.Item(22, 211).Value = 16
$ws.Cells.Item(22, 212).Value = 8
$ws.Cells.Item(22, 213).Value = 9
$ws.Cells.Item(23, 211).Value = 1
$ws.Cells.Item(23, 212).Value = 4
$ws.Cells.Item(23, 213).Value = 2
$ws.Cells.Item(24, 211).Value = 21
$ws.Cells.Item(24, 212).Value = 26
$ws.Cells.Item(24, 213).Value = 18
$ws.Cells.Item(25, 211).Value = 19
$ws.Cells.Item(25, 212).Value = 53.8
$ws.Cells.Item(25, 213).Value = 38.9
$ws.Cells.Item(26, 211).Value = 74.5
$ws.Cells.Item(26, 212).Value = 20.71
$ws.Cells.Item(26, 213).Value = 37.86
$ws.Cells.Item(27, 211).Value = 14.19
$ws.Cells.Item(27, 212).Value = 11.15
$ws.Cells.Item(27, 213).Value = 14.72
$ws.Cells.Item(28, 211).Value = 28
$ws.Cells.Item(28, 212).Value = 36
$ws.Cells.Item(28, 213).Value = 26
$ws.Cells.Item(29, 211).Value = 34
$ws.Cells.Item(29, 212).Value = 36
$ws.Cells.Item(29, 213).Value = 42
$ws.Cells.Item(30, 211).Value = 31
$ws.Cells.Item(30, 212).Value = 36
$ws.Cells.Item(30, 213).Value = 27
$ws.Cells.Item(31, 211).Value = 40
$ws.Cells.Item(31, 212).Value = 53
$ws.Cells.Item(31, 213).Value = 47
$ws.Cells.Item(32, 211).Value = 1.9
$ws.Cells.Item(32, 212).Value = 2.04
$ws.Cells.Item(32, 213).Value = 2.61
$ws.Cells.Item(33, 211).Value = 10
$ws.Cells.Item(33, 212).Value = 3.79
$ws.Cells.Item(33, 213).Value = 6.71
$ws.Cells.Item(34, 211).Value = 50
$ws.Cells.Item(34, 212).Value = 41.5
$ws.Cells.Item(34, 213).Value = 34
$ws.Cells.Item(35, 211).Value = 10
$ws.Cells.Item(35, 212).Value = 26.4
$ws.Cells.Item(35, 213).Value = 14.9
$ws.Cells.Item(36, 211).Value = 188.8
$ws.Cells.Item(36, 212).Value = 188.3
$ws.Cells.Item(36, 213).Value = 188.6
$ws.Cells.Item(37, 211).Value = 89.2
$ws.Cells.Item(37, 212).Value = 88.5
$ws.Cells.Item(37, 213).Value = 88
$ws.Cells.Item(38, 211).Value = 24.24
$ws.Cells.Item(38, 212).Value = 25.74
$ws.Cells.Item(38, 213).Value = 25.49
$ws.Cells.Item(39, 211).Value = 74.5
$ws.Cells.Item(39, 212).Value = 99.3
$ws.Cells.Item(39, 213).Value = 96.5
$ws.Cells.Item(40, 211).Value = 8
$ws.Cells.Item(40, 212).Value = 6
$ws.Cells.Item(40, 213).Value = 6
$ws.Cells.Item(41, 211).Value = 7
$ws.Cells.Item(41, 212).Value = 6
$ws.Cells.Item(41, 213).Value = 6
$ws.Cells.Item(42, 211).Value = 4
$ws.Cells.Item(42, 212).Value = 5
$ws.Cells.Item(42, 213).Value = 6
$ws.Cells.Item(43, 211).Value = 3
$ws.Cells.Item(43, 212).Value = 5
$ws.Cells.Item(43, 213).Value = 4
$ws.Cells.Item(44, 211).Value = 97
$ws.Cells.Item(44, 212).Value = 113
$ws.Cells.Item(44, 213).Value = 114
$ws.Cells.Item(45, 211).Value = 199
$ws.Cells.Item(45, 212).Value = 171
$ws.Cells.Item(45, 213).Value = 151
$ws.Cells.Item(46, 211).Value = 217
$ws.Cells.Item(46, 212).Value = 211
$ws.Cells.Item(46, 213).Value = 182
$ws.Cells.Item(47, 211).Value = 72.8
$ws.Cells.Item(47, 212).Value = 72.8
$ws.Cells.Item(47, 213).Value = 68.7
$ws.Cells.Item(48, 211).Value = 34
$ws.Cells.Item(48, 212).Value = 36
$ws.Cells.Item(48, 213).Value = 42
$ws.Cells.Item(49, 211).Value = 8
$ws.Cells.Item(49, 212).Value = 6
$ws.Cells.Item(49, 213).Value = 12
$ws.Cells.Item(50, 211).Value = 8
$ws.Cells.Item(50, 212).Value = 10
$ws.Cells.Item(50, 213).Value = 10
$ws.Cells.Item(51, 211).Value = 28
$ws.Cells.Item(51, 212).Value = 36
$ws.Cells.Item(51, 213).Value = 26
$ws.Cells.Item(52, 211).Value = 31
$ws.Cells.Item(52, 212).Value = 36
$ws.Cells.Item(52, 213).Value = 27
$ws.Cells.Item(53, 211).Value = 33
$ws.Cells.Item(53, 212).Value = 41
$ws.Cells.Item(53, 213).Value = 32
$ws.Cells.Item(54, 211).Value = 1
$ws.Cells.Item(54, 212).Value = 1
$ws.Cells.Item(54, 213).Value = 1
$ws.Cells.Item(55, 211).Value = 1
$ws.Cells.Item(55, 212).Value = 9
$ws.Cells.Item(55, 213).Value = 5
$ws.Cells.Item(56, 211).Value = 25
$ws.Cells.Item(56, 212).Value = 64.3
$ws.Cells.Item(56, 213).Value = 71.40000000000001
$ws.Cells.Item(57, 211).Value = 183
$ws.Cells.Item(57, 212).Value = 159
$ws.Cells.Item(57, 213).Value = 171
$ws.Cells.Item(58, 211).Value = 112
$ws.Cells.Item(58, 212).Value = 150
$ws.Cells.Item(58, 213).Value = 137
$ws.Cells.Item(59, 211).Value = 295
$ws.Cells.Item(59, 212).Value = 309
$ws.Cells.Item(59, 213).Value = 308
$ws.Cells.Item(60, 211).Value = 1.63
$ws.Cells.Item(60, 212).Value = 1.06
$ws.Cells.Item(60, 213).Value = 1.25
$ws.Cells.Item(61, 211).Value = 80
$ws.Cells.Item(61, 212).Value = 64
$ws.Cells.Item(61, 213).Value = 77
$ws.Cells.Item(62, 211).Value = 39
$ws.Cells.Item(62, 212).Value = 47
$ws.Cells.Item(62, 213).Value = 50
$ws.Cells.Item(63, 211).Value = 23
$ws.Cells.Item(63, 212).Value = 23
$ws.Cells.Item(63, 213).Value = 25
$ws.Cells.Item(64, 211).Value = 11
$ws.Cells.Item(64, 212).Value = 12
$ws.Cells.Item(64, 213).Value = 19
$ws.Cells.Item(65, 211).Value = 18
$ws.Cells.Item(65, 212).Value = 9
$ws.Cells.Item(65, 213).Value = 13
$ws.Cells.Item(66, 211).Value = 12
$ws.Cells.Item(66, 212).Value = 11
$ws.Cells.Item(66, 213).Value = 8
$ws.Cells.Item(67, 211).Value = 9
$ws.Cells.Item(67, 212).Value = 8
$ws.Cells.Item(67, 213).Value = 6
$ws.Cells.Item(68, 211).Value = 7
$ws.Cells.Item(68, 212).Value = 4
$ws.Cells.Item(68, 213).Value = 4
$ws.Cells.Item(69, 211).Value = 3
$ws.Cells.Item(69, 212).Value = 2
$ws.Cells.Item(69, 213).Value = 0
$ws.Cells.Item(70, 211).Value = 22
$ws.Cells.Item(70, 212).Value = 17
$ws.Cells.Item(70, 213).Value = 12
$ws.Cells.Item(71, 211).Value = 54.5
$ws.Cells.Item(71, 212).Value = 64.7
$ws.Cells.Item(71, 213).Value = 66.7
$ws.Cells.Item(72, 211).Value = 24.58
$ws.Cells.Item(72, 212).Value = 28.09
$ws.Cells.Item(72, 213).Value = 38.5
$ws.Cells.Item(73, 211).Value = 13.41
$ws.Cells.Item(73, 212).Value = 18.18
$ws.Cells.Item(73, 213).Value = 25.67
$ws.Cells.Item(74, 211).Value = 20
$ws.Cells.Item(74, 212).Value = 27
$ws.Cells.Item(74, 213).Value = 32
$ws.Cells.Item(75, 211).Value = 57
$ws.Cells.Item(75, 212).Value = 43
$ws.Cells.Item(75, 213).Value = 44
$ws.Cells.Item(76, 211).Value = 35
$ws.Cells.Item(76, 212).Value = 38
$ws.Cells.Item(76, 213).Value = 39
$ws.Cells.Item(77, 211).Value = 44
$ws.Cells.Item(77, 212).Value = 47
$ws.Cells.Item(77, 213).Value = 35
$ws.Cells.Item(78, 211).Value = 2
$ws.Cells.Item(78, 212).Value = 2.76
$ws.Cells.Item(78, 213).Value = 2.92
$ws.Cells.Item(79, 211).Value = 3.67
$ws.Cells.Item(79, 212).Value = 4.27
$ws.Cells.Item(79, 213).Value = 4.38
$ws.Cells.Item(80, 211).Value = 43.2
$ws.Cells.Item(80, 212).Value = 31.9
$ws.Cells.Item(80, 213).Value = 34.3
$ws.Cells.Item(81, 211).Value = 27.3
$ws.Cells.Item(81, 212).Value = 23.4
$ws.Cells.Item(81, 213).Value = 22.9
$ws.Cells.Item(82, 211).Value = 186.1
$ws.Cells.Item(82, 212).Value = 187.7
$ws.Cells.Item(82, 213).Value = 189.6
$ws.Cells.Item(83, 211).Value = 85.59999999999999
$ws.Cells.Item(83, 212).Value = 85.5
$ws.Cells.Item(83, 213).Value = 86.40000000000001
$ws.Cells.Item(84, 211).Value = 25.24
$ws.Cells.Item(84, 212).Value = 24.8
$ws.Cells.Item(84, 213).Value = 24.49
$ws.Cells.Item(85, 211).Value = 85.8
$ws.Cells.Item(85, 212).Value = 82.2
$ws.Cells.Item(85, 213).Value = 68.09999999999999
$ws.Cells.Item(86, 211).Value = 11
$ws.Cells.Item(86, 212).Value = 7
$ws.Cells.Item(86, 213).Value = 12
$ws.Cells.Item(87, 211).Value = 4
$ws.Cells.Item(87, 212).Value = 7
$ws.Cells.Item(87, 213).Value = 5
$ws.Cells.Item(88, 211).Value = 1
$ws.Cells.Item(88, 212).Value = 6
$ws.Cells.Item(88, 213).Value = 2
$ws.Cells.Item(89, 211).Value = 6
$ws.Cells.Item(89, 212).Value = 2
$ws.Cells.Item(89, 213).Value = 3
$ws.Cells.Item(90, 211).Value = 103
$ws.Cells.Item(90, 212).Value = 120
$ws.Cells.Item(90, 213).Value = 116
$ws.Cells.Item(91, 211).Value = 182
$ws.Cells.Item(91, 212).Value = 189
$ws.Cells.Item(91, 213).Value = 182
$ws.Cells.Item(92, 211).Value = 225
$ws.Cells.Item(92, 212).Value = 228
$ws.Cells.Item(92, 213).Value = 223
$ws.Cells.Item(93, 211).Value = 76.3
$ws.Cells.Item(93, 212).Value = 73.8
$ws.Cells.Item(93, 213).Value = 72.40000000000001
$ws.Cells.Item(94, 211).Value = 57
$ws.Cells.Item(94, 212).Value = 43
$ws.Cells.Item(94, 213).Value = 44
$ws.Cells.Item(95, 211).Value = 13
$ws.Cells.Item(95, 212).Value = 8
$ws.Cells.Item(95, 213).Value = 15
$ws.Cells.Item(96, 211).Value = 12
$ws.Cells.Item(96, 212).Value = 10
$ws.Cells.Item(96, 213).Value = 6
$ws.Cells.Item(97, 211).Value = 20
$ws.Cells.Item(97, 212).Value = 27
$ws.Cells.Item(97, 213).Value = 32
$ws.Cells.Item(98, 211).Value = 35
$ws.Cells.Item(98, 212).Value = 38
$ws.Cells.Item(98, 213).Value = 39
$ws.Cells.Item(99, 211).Value = 24
$ws.Cells.Item(99, 212).Value = 44
$ws.Cells.Item(99, 213).Value = 30
$ws.Cells.Item(100, 211).Value = 5
$ws.Cells.Item(100, 212).Value = 3
$ws.Cells.Item(100, 213).Value = 10
$ws.Cells.Item(101, 211).Value = 9
$ws.Cells.Item(101, 212).Value = 8
$ws.Cells.Item(101, 213).Value = 6
$ws.Cells.Item(102, 211).Value = 75
$ws.Cells.Item(102, 212).Value = 72.7
$ws.Cells.Item(102, 213).Value = 75
